# Auto-generated Excel COM-interop script to apply the scheduled market-price data refresh
# described by the workbook diff (static currentAveragePrice / LevePrice / LeveProfit columns H:N).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 10773854
$ws.Range("I19").Value = 9391917
$ws.Range("J19").Value = 12501276
$ws.Range("K19").Value = 9391917
$ws.Range("L19").Value = 12501276
$ws.Range("M19").Value = -9391742
$ws.Range("N19").Value = -12501626
$ws.Range("H21").Value = 23394.363
$ws.Range("I21").Value = 70019
$ws.Range("J21").Value = 13033.333
$ws.Range("K21").Value = 70019
$ws.Range("L21").Value = 13033.333
$ws.Range("M21").Value = -69551
$ws.Range("N21").Value = -13969.333
$ws.Range("H23").Value = 23394.363
$ws.Range("I23").Value = 70019
$ws.Range("J23").Value = 13033.333
$ws.Range("K23").Value = 70019
$ws.Range("L23").Value = 13033.333
$ws.Range("M23").Value = -69785
$ws.Range("N23").Value = -13501.333
$ws.Range("H28").Value = 1039
$ws.Range("I28").Value = 329.23077
$ws.Range("J28").Value = 2357.1428
$ws.Range("K28").Value = 329.23077
$ws.Range("L28").Value = 2357.1428
$ws.Range("M28").Value = 155.76923
$ws.Range("N28").Value = -3327.1428
$ws.Range("H33").Value = 189
$ws.Range("I33").Value = 40.307693
$ws.Range("J33").Value = 833.3333
$ws.Range("K33").Value = 40.307693
$ws.Range("L33").Value = 833.3333
$ws.Range("M33").Value = 188.692307
$ws.Range("N33").Value = -1291.3333
$ws.Range("H40").Value = 1864.3529
$ws.Range("I40").Value = 1840.1333
$ws.Range("J40").Value = 2046
$ws.Range("K40").Value = 1840.1333
$ws.Range("L40").Value = 2046
$ws.Range("M40").Value = -1665.1333
$ws.Range("N40").Value = -2396
$ws.Range("H64").Value = 3877.3547
$ws.Range("I64").Value = 3768.75
$ws.Range("J64").Value = 3993.2
$ws.Range("K64").Value = 3768.75
$ws.Range("L64").Value = 3993.2
$ws.Range("M64").Value = -3520.75
$ws.Range("N64").Value = -4489.2
$ws.Range("H67").Value = 3877.3547
$ws.Range("I67").Value = 3768.75
$ws.Range("J67").Value = 3993.2
$ws.Range("K67").Value = 3768.75
$ws.Range("L67").Value = 3993.2
$ws.Range("M67").Value = -2910.75
$ws.Range("N67").Value = -5709.2
$ws.Range("H96").Value = 1500
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1500
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 4500
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -7246
$ws.Range("H99").Value = 335.7
$ws.Range("I99").Value = 357
$ws.Range("J99").Value = 286
$ws.Range("K99").Value = 1071
$ws.Range("L99").Value = 858
$ws.Range("M99").Value = 427
$ws.Range("N99").Value = -3854
$ws.Range("H100").Value = 1445.6875
$ws.Range("I100").Value = 1471.2222
$ws.Range("J100").Value = 1412.8572
$ws.Range("K100").Value = 1471.2222
$ws.Range("L100").Value = 1412.8572
$ws.Range("M100").Value = -930.2221999999999
$ws.Range("N100").Value = -2494.8572
$ws.Range("H125").Value = 4346.6665
$ws.Range("J125").Value = 4346.6665
$ws.Range("L125").Value = 39119.9985
$ws.Range("N125").Value = -44039.9985
$ws.Range("H141").Value = 857
$ws.Range("I141").Value = 746.25
$ws.Range("J141").Value = 1300
$ws.Range("K141").Value = 2238.75
$ws.Range("L141").Value = 3900
$ws.Range("M141").Value = 2941.25
$ws.Range("N141").Value = -14260

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1179.091
$ws.Range("I102").Value = 1110
$ws.Range("K102").Value = 1110
$ws.Range("M102").Value = 512
$ws.Range("H122").Value = 2120
$ws.Range("I122").Value = 2034.1428
$ws.Range("J122").Value = 2240.2
$ws.Range("K122").Value = 6102.428400000001
$ws.Range("L122").Value = 6720.599999999999
$ws.Range("M122").Value = -3652.428400000001
$ws.Range("N122").Value = -11620.6
$ws.Range("H132").Value = 2253.6365
$ws.Range("I132").Value = 1239
$ws.Range("J132").Value = 3471.2
$ws.Range("K132").Value = 3717
$ws.Range("L132").Value = 10413.6
$ws.Range("M132").Value = -1187
$ws.Range("N132").Value = -15473.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 437.66666
$ws.Range("I94").Value = 461.2
$ws.Range("J94").Value = 320
$ws.Range("K94").Value = 461.2
$ws.Range("L94").Value = 320
$ws.Range("M94").Value = -10.19999999999999
$ws.Range("N94").Value = -1222
$ws.Range("H105").Value = 2699.63
$ws.Range("I105").Value = 1648.8889
$ws.Range("J105").Value = 2803.5496
$ws.Range("K105").Value = 1648.8889
$ws.Range("L105").Value = 2803.5496
$ws.Range("M105").Value = 98.11110000000008
$ws.Range("N105").Value = -6297.5496
$ws.Range("H134").Value = 2409.35
$ws.Range("I134").Value = 2231.8333
$ws.Range("K134").Value = 6695.499899999999
$ws.Range("M134").Value = -4160.499899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 788.9286
$ws.Range("I34").Value = 116.666664
$ws.Range("J34").Value = 972.2727
$ws.Range("K34").Value = 349.999992
$ws.Range("L34").Value = 2916.8181
$ws.Range("M34").Value = -265.999992
$ws.Range("N34").Value = -3084.8181
$ws.Range("H39").Value = 3166.3333
$ws.Range("J39").Value = 3166.3333
$ws.Range("L39").Value = 9498.999899999999
$ws.Range("N39").Value = -10086.9999
$ws.Range("H55").Value = 3706.4167
$ws.Range("J55").Value = 4018.4546
$ws.Range("L55").Value = 12055.3638
$ws.Range("N55").Value = -12409.3638
$ws.Range("H68").Value = 310
$ws.Range("I68").Value = 310
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 930
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -119
$ws.Range("N68").ClearContents()
$ws.Range("H69").Value = 2167.1052
$ws.Range("I69").Value = 1435.8334
$ws.Range("J69").Value = 2504.6155
$ws.Range("K69").Value = 4307.5002
$ws.Range("L69").Value = 7513.8465
$ws.Range("M69").Value = -3496.5002
$ws.Range("N69").Value = -9135.8465
$ws.Range("H71").Value = 310
$ws.Range("I71").Value = 310
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 2790
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 1266
$ws.Range("N71").ClearContents()
$ws.Range("H72").Value = 2167.1052
$ws.Range("I72").Value = 1435.8334
$ws.Range("J72").Value = 2504.6155
$ws.Range("K72").Value = 12922.5006
$ws.Range("L72").Value = 22541.5395
$ws.Range("M72").Value = -8866.500599999999
$ws.Range("N72").Value = -30653.5395
$ws.Range("H82").Value = 4806
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 4806
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 14418
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -15230
$ws.Range("H85").Value = 4806
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 4806
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 14418
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -17226

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1845.1
$ws.Range("I97").Value = 1916.6666
$ws.Range("J97").Value = 1814.4286
$ws.Range("K97").Value = 1916.6666
$ws.Range("L97").Value = 1814.4286
$ws.Range("M97").Value = -1420.6666
$ws.Range("N97").Value = -2806.4286
$ws.Range("H122").Value = 20001690
$ws.Range("I122").Value = 25001362
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 75004086
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -75001636
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2313
$ws.Range("I132").Value = 2268.9524
$ws.Range("J132").Value = 2384.1538
$ws.Range("K132").Value = 6806.8572
$ws.Range("L132").Value = 7152.4614
$ws.Range("M132").Value = -4276.8572
$ws.Range("N132").Value = -12212.4614

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 11628.8
$ws.Range("I93").Value = 18881.666
$ws.Range("K93").Value = 18881.666
$ws.Range("M93").Value = -17633.666
$ws.Range("H133").Value = 50707.418
$ws.Range("J133").Value = 50707.418
$ws.Range("L133").Value = 50707.418
$ws.Range("N133").Value = -55767.418

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 654.1
$ws.Range("I100").Value = 345.36365
$ws.Range("J100").Value = 1031.4445
$ws.Range("K100").Value = 690.7273
$ws.Range("L100").Value = 2062.889
$ws.Range("M100").Value = -149.7273
$ws.Range("N100").Value = -3144.889
$ws.Range("H139").Value = 37877.57
$ws.Range("J139").Value = 37877.57
$ws.Range("L139").Value = 37877.57
$ws.Range("N139").Value = -48157.57
